$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new value for cell B4 (indicator title)
$ws.Range("B4").Value = " 3.5.2.1 Incidence of alcohol dependence per 100 000 people"

# Update the selection / view on the sheet (no longer scrolled to A21, selection now C6)
$ws.Range("C6").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
